# CU-1vhwepw added ntc for api Register
# Adds two new negative/positive test-data rows to the "apiTest" sheet:
#   row 3: invalid email "srdjan.rados@htecgroup" (missing ".com"), same
#          password/first/last name as the existing row.
#   row 4: the original valid email "srdjan.rados@htecgroup.com" again,
#          paired with a different password "Qwertysha".
# Both new A-column values get a mailto: hyperlink, matching the existing
# row 2 convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiTest")

# ---- Row 3 ----
$ws.Range("A3").Value = "srdjan.rados@htecgroup"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:srdjan.rados@htecgroup", "", "", "srdjan.rados@htecgroup") | Out-Null
# Hyperlinks.Add auto-applies Excel's blue/underline hyperlink style; the
# source workbook keeps these cells on the plain default style, so strip
# the styling back off again.
$ws.Range("A3").Font.Underline = $false
$ws.Range("A3").Font.ColorIndex = 1
$ws.Range("B3").Value = "Qwertysha1@"
$ws.Range("C3").Value = "Srdjan"
$ws.Range("D3").Value = "Rados"

# ---- Row 4 ----
$ws.Range("A4").Value = "srdjan.rados@htecgroup.com"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:srdjan.rados@htecgroup.com", "", "", "srdjan.rados@htecgroup.com") | Out-Null
$ws.Range("A4").Font.Underline = $false
$ws.Range("A4").Font.ColorIndex = 1
$ws.Range("B4").Value = "Qwertysha"
$ws.Range("C4").Value = "Srdjan"
$ws.Range("D4").Value = "Rados"

# Widen column A (longer hyperlink text) and nudge column D back in,
# matching the resize that accompanied the new rows.
$ws.Columns.Item(1).ColumnWidth = 35.8
$ws.Columns.Item(4).ColumnWidth = 24.5

# Selection moves to the newly added last row.
$ws.Range("A4").Select() | Out-Null
